$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Add new cell G18 = "questionnary" (shared string index 28)
$ws.Range("G18").Value = "questionnary"

# Add new cell E19 = "arrayOfGoodAnswers" (new shared string)
$ws.Range("E19").Value = "arrayOfGoodAnswers"

# Update selection / view: clear the frozen/scrolled topLeftCell and move selection to E20
$ws.Range("E20").Select()

$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
